$d = $word.ActiveDocument

function Replace-AllText($findText, $replaceText) {
    $searchStart = 0
    while ($true) {
        $docEnd = $d.Content.End
        if ($searchStart -ge $docEnd) { break }
        $rng = $d.Range($searchStart, $docEnd)
        $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if (-not $found) { break }
        $rng.Text = $replaceText
        $searchStart = $rng.End
    }
}

Replace-AllText "Wafungwa na peremende - manukuu:" "Prisoners and candies - subtitles:"
Replace-AllText "**mazungumzo huanza saa 55 ya pili sio 27 kwa sababu ya klipu ya utangulizi. Nilirekebisha nyakati ipasavyo. -John Argentino" "**dialogue starts at second 55 not 27 because of the intro clip. I adjusted the times accordingly. -John Argentino"
Replace-AllText "[Muziki]" "[Music]"
Replace-AllText "wanahisabati wanne mkali wanachukuliwa" "four bright mathematicians are taken into"
Replace-AllText "chini ya ulinzi na kuwekwa gerezani kwa sababu walijaribu" "custody and put in jail because they tried"
Replace-AllText "kumshawishi mwanamke mzee kuwa Goedel's" "to convince an old lady that the Goedel's"
Replace-AllText "nadharia za kutokamilika ni kweli. Kila" "incompleteness theorems are true. Every"
Replace-AllText "mtaalamu wa hisabati ana kiini chake ambacho sisi" "mathematician has his own cell that we"
Replace-AllText "inaweza kuhesabu na nambari kutoka 1 hadi 4." "can enumerate with a number from 1 to 4."
Replace-AllText "kabla ya kuingia kwenye seli fulani" "before entering the cell a certain"
Replace-AllText "idadi ya peremende kubwa kuliko " "number of candies greater than "
Replace-AllText "e: AU " "e: OR "
Replace-AllText "SAWA NA" "EQUAL TO"
Replace-AllText " 1 ni" " 1 is"
Replace-AllText "wanapewa kila mtaalamu wa hisabati na wao" "given to every mathematician and they"
Replace-AllText "wanaambiwa wana peremende 11 kwa jumla." "are told they have 11 candies in total."
Replace-AllText "lakini kila mtu anajua idadi yake tu" "but everyone knows only his number of"
Replace-AllText "pipi na jumla. 1 na sio" "candies and the total. 1 and is not"
Replace-AllText "kuruhusiwa kuuliza nambari zingine." "allowed to ask for the others number."
Replace-AllText "kisha mwanahisabati wa kwanza anauliza" "then the first mathematician asks the"
Replace-AllText "pili: 'namba 2 unajua kama wewe" "second: 'number 2 do you know if you"
Replace-AllText "kuwa na peremende nyingi kuliko mimi?' ya pili" "have more candies than me?' the second"
Replace-AllText "mwanahisabati anajibu hana. Kisha yeye" "mathematician answers he doesn't. Then he"
Replace-AllText "anauliza kwa nambari 3: 'unajua kama unayo" "asks to number 3: 'do you know if you have"
Replace-AllText "pipi zaidi kuliko mimi?'" "more candy than me?'"
Replace-AllText "mwanahisabati wa tatu anajibu: 'hapana niko" "the third mathematician answers: 'no I'm"
Replace-AllText "samahani sifanyi'. Katika hatua hii ya nne" "sorry I don't'. At this point the fourth"
Replace-AllText "mtaalamu wa hisabati anasema: 'jamani mnafahamu" "mathematician says: 'hey guys you know"
Replace-AllText "nini, najua hasa pipi ngapi" "what, I know exactly how many candies"
Replace-AllText "kila mtu ana hapa'. Cha kushangaza hata" "everyone has here'. Surprisingly even the"
Replace-AllText "wanahisabati wengine watatu wanasema hivyo sasa" "other three mathematicians say that now"
Replace-AllText "wanajua kila mtu ana pipi ngapi" "they know how many candies everyone has"
Replace-AllText "kwa hivyo swali ni: unaweza kujua" "so the question is: can you figure out"
Replace-AllText "idadi ya pipi kila mfungwa ana" "the number of candies every prisoner has"

# Remove the stray single-space run that sits between the ")" run and the
# " 1 is" run (the diff drops this whole <w:r> entirely).
$rng = $d.Content
$found = $rng.Find.Execute(")", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $spaceRng = $d.Range($rng.End, $rng.End + 1)
    if ($spaceRng.Text -eq " ") {
        $spaceRng.Delete()
    }
}
